$wb = $excel.ActiveWorkbook

# Rename the "BCR" worksheet to "BCRbQ". Excel will automatically update any
# formulas elsewhere in the workbook that reference the sheet by name.
$bcrSheet = $wb.Worksheets.Item("BCR")
$bcrSheet.Name = "BCRbQ"

# Update the title cell on the "About" sheet to match the new sheet name.
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A1").Value = "BCRbQ BAU Capacity Retirements before Quantization"
